# Updated cryptos list on Wed Apr 26 13:42:48 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.736.01"

$ws.Range("D3").Value = "1.949.64"
$ws.Range("E3").Value = "  +6.73%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "342.54"
$ws.Range("E5").Value = "  +3.48%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "0.4793"
$ws.Range("E7").Value = "  +4.68%  "

$ws.Range("D8").Value = "0.4142"
$ws.Range("E8").Value = "  +8.47%  "

$ws.Range("D9").Value = "'48.80"
$ws.Range("E9").Value = "  +5.27%  "

$ws.Range("D10").Value = "0.08281"
$ws.Range("E10").Value = "  +4.99%  "

$ws.Range("E11").Value = "  +8.58%  "

$ws.Range("D12").Value = "22.74"
$ws.Range("E12").Value = "  +7.95%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.937.38"
$ws.Range("E13").Value = "  +5.94%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.177"
$ws.Range("E14").Value = "  +5.15%  "

$ws.Range("D15").Value = "7.447"
$ws.Range("E15").Value = "  +5.05%  "

$ws.Range("E16").Value = "  +3.58%  "

$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("E18").Value = "  +4.22%  "

$ws.Range("D19").Value = "0.06722"
$ws.Range("E19").Value = "  +1.76%  "

$ws.Range("D20").Value = "18.05"
$ws.Range("E20").Value = "  +5.03%  "

$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").Value = "29.704.73"
$ws.Range("E22").Value = "  +8.36%  "

$ws.Range("D23").Value = "5.632"
$ws.Range("E23").Value = "  +5.96%  "

$ws.Range("E24").Value = "  +4.13%  "

$ws.Range("D25").Value = "2.277"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "2.166.68"
$ws.Range("E26").Value = "  +5.93%  "

$ws.Range("E27").Value = "  +3.62%  "

$ws.Range("D28").Value = "20.15"
$ws.Range("E28").Value = "  +3.99%  "

$ws.Range("D29").Value = "2.201"
$ws.Range("E29").Value = "  +6.35%  "

$ws.Range("D30").Value = "5.647"
$ws.Range("E30").Value = "  +6.86%  "

$ws.Range("D31").Value = "122.61"

$ws.Range("D32").Value = "1.026"
$ws.Range("E32").Value = "  +9.22%  "

$ws.Range("E33").Value = "  +3.77%  "

$ws.Range("D34").Value = "1.475"
$ws.Range("E34").Value = "  +11.28%  "

$ws.Range("E35").Value = "  +2.97%  "

$ws.Range("D36").Value = "5.495"
$ws.Range("E36").Value = "  +5.07%  "

$ws.Range("D37").Value = "0.02317"
$ws.Range("E37").Value = "  +6.19%  "

$ws.Range("D38").Value = "0.06246"
$ws.Range("E38").Value = "  +5.34%  "

$ws.Range("D39").Value = "8.717"
$ws.Range("E39").Value = "  +7.11%  "

$ws.Range("D40").Value = "1.199"
$ws.Range("E40").Value = "  +4.79%  "

$ws.Range("D41").Value = "0.6113"
$ws.Range("E41").Value = "  +5.89%  "

$ws.Range("D42").Value = "10.76"
$ws.Range("E42").Value = "  +7.71%  "

$ws.Range("D43").Value = "0.1909"
$ws.Range("E43").Value = "  +4.56%  "

$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").Value = "1.282"
$ws.Range("E45").Value = "  +1.46%  "

$ws.Range("D46").Value = "12.65"
$ws.Range("E46").Value = "  +5.82%  "

$ws.Range("D47").Value = "0.5729"
$ws.Range("E47").Value = "  +5.46%  "

$ws.Range("D48").Value = "2.326"
$ws.Range("E48").Value = "  +27.41%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.07476"
$ws.Range("E49").Value = "  +13.77%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "2.001"
$ws.Range("E50").Value = "  +6.91%  "

$ws.Range("D51").Value = "114.23"
$ws.Range("E51").Value = "  +4.16%  "
